# Applies the "add 2022-Q3 data" edit to the 688510-航亚科技 workbook.
#
# Summary of the change (per the commit's canonical-OOXML diff):
#  1. A brand-new worksheet "2022-Q3" is inserted right before the existing
#     "2022-Q2" worksheet, holding new per-fund position data for 001643 /
#     001644 (everything else about that sheet - headers, styles, fund
#     codes/names - mirrors the existing quarter sheets).
#  2. The "总计" (summary) sheet gains one new row at the top (row 2) for
#     2022-Q3, pushing the previously-existing rows down by one; every
#     other row keeps exactly the data it already had.
#  All other sheets/tabs keep their name and contents unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" sheet by duplicating "2022-Q2" (so it
# inherits the same headers / column styles / number formats), then
# overwrite only the data cells that actually differ.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)                      # new copy is placed immediately before $q2
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2: fund 001643 / 汇丰晋信智造先锋股票A (code/name copied as-is)
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "15.68"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "94.47"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.99"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.6256"
$q3.Range("H2").Value = 6

# Row 3: fund 001644 / 汇丰晋信智造先锋股票C (code/name copied as-is)
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "8.27"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "94.47"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "3.99"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.3300"
$q3.Range("H3").Value = 6

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q3 row to the "总计" summary sheet and re-write
# the (unchanged) data for every other quarter so the table reads, top
# to bottom: 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$rows = @(
    @(0, "2022-Q3", 2, 0.96),
    @(1, "2022-Q2", 2, 1.57),
    @(2, "2022-Q1", 2, 1.5),
    @(3, "2021-Q4", 2, 0.83),
    @(4, "2021-Q3", 2, 1.65),
    @(5, "2021-Q2", 2, 1.59),
    @(6, "2021-Q1", 6, 1.44)
)

$r = 2
foreach ($row in $rows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
